$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1, columns E..BL previously held the year labels as text
# ("1960 [YR1960]" ... "2019 [YR2019]"). Replace them with the actual
# numeric year values and left-align them (matches the new cellXfs style).
$ws.Range("E1:BL1").HorizontalAlignment = -4131

for ($i = 0; $i -lt 60; $i++) {
    $ws.Cells.Item(1, 5 + $i).Value = 1960 + $i
}

# Update the selection to reflect the edited range.
[void]$ws.Range("E1:BL1").Select()
